$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value2 = 13253.154
$ws.Range("J51").Value2 = 8072
$ws.Range("L51").Value2 = 8072
$ws.Range("N51").Value2 = -9040

$ws.Range("H70").Value2 = 21609570
$ws.Range("J70").Value2 = 23812052
$ws.Range("L70").Value2 = 71436156
$ws.Range("N70").Value2 = -71436696

$ws.Range("H73").Value2 = 21609570
$ws.Range("J73").Value2 = 23812052
$ws.Range("L73").Value2 = 71436156
$ws.Range("N73").Value2 = -71438028

$ws.Range("H74").Value2 = 115394250
$ws.Range("I74").Value2 = 214288740
$ws.Range("K74").Value2 = 214288740
$ws.Range("M74").Value2 = -214287804

$ws.Range("H76").Value2 = 33335834
$ws.Range("I76").Value2 = 100000000
$ws.Range("K76").Value2 = 100000000
$ws.Range("M76").Value2 = -99999685

$ws.Range("H77").Value2 = 115394250
$ws.Range("I77").Value2 = 214288740
$ws.Range("K77").Value2 = 1071443700
$ws.Range("M77").Value2 = -1071439020

$ws.Range("H79").Value2 = 33335834
$ws.Range("I79").Value2 = 100000000
$ws.Range("K79").Value2 = 100000000
$ws.Range("M79").Value2 = -99998908

$ws.Range("H137").Value2 = 2965
$ws.Range("I137").Value2 = 2754.8215
$ws.Range("K137").Value2 = 8264.4645
$ws.Range("M137").Value2 = -5714.4645

$ws.Range("H138").Value2 = 4972.696
$ws.Range("J138").Value2 = 8076.45
$ws.Range("L138").Value2 = 24229.35
$ws.Range("N138").Value2 = -34509.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 1281201.8
$ws.Range("I32").Value2 = 1348631.9
$ws.Range("J32").Value2 = 27002.2
$ws.Range("K32").Value2 = 1348631.9
$ws.Range("L32").Value2 = 27002.2
$ws.Range("M32").Value2 = -1348344.9
$ws.Range("N32").Value2 = -27576.2

$ws.Range("H45").Value2 = 6310.353
$ws.Range("I45").Value2 = 5250.923
$ws.Range("J45").Value2 = 9753.5
$ws.Range("K45").Value2 = 5250.923
$ws.Range("L45").Value2 = 9753.5
$ws.Range("M45").Value2 = -4873.923
$ws.Range("N45").Value2 = -10507.5

$ws.Range("H63").Value2 = 1541.625
$ws.Range("I63").Value2 = 1111
$ws.Range("J63").Value2 = 1800
$ws.Range("K63").Value2 = 1111
$ws.Range("L63").Value2 = 1800
$ws.Range("M63").Value2 = -425
$ws.Range("N63").Value2 = -3172

$ws.Range("H66").Value2 = 1541.625
$ws.Range("I66").Value2 = 1111
$ws.Range("J66").Value2 = 1800
$ws.Range("K66").Value2 = 5555
$ws.Range("L66").Value2 = 9000
$ws.Range("M66").Value2 = -2123
$ws.Range("N66").Value2 = -15864

$ws.Range("H122").Value2 = 2476.8108
$ws.Range("I122").Value2 = 1699.2142
$ws.Range("J122").Value2 = 4896
$ws.Range("K122").Value2 = 5097.642599999999
$ws.Range("L122").Value2 = 14688
$ws.Range("M122").Value2 = -2647.642599999999
$ws.Range("N122").Value2 = -19588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 41669676
$ws.Range("I107").Value2 = 56251564
$ws.Range("K107").Value2 = 56251564
$ws.Range("M107").Value2 = -56249644

$ws.Range("H134").Value2 = 5636.6274
$ws.Range("I134").Value2 = 2126.1538
$ws.Range("K134").Value2 = 6378.4614
$ws.Range("M134").Value2 = -3843.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 6165.1846
$ws.Range("I31").Value2 = 2922.0488
$ws.Range("J31").Value2 = 11705.542
$ws.Range("K31").Value2 = 2922.0488
$ws.Range("L31").Value2 = 11705.542
$ws.Range("M31").Value2 = -2627.0488
$ws.Range("N31").Value2 = -12295.542

$ws.Range("H34").Value2 = 6165.1846
$ws.Range("I34").Value2 = 2922.0488
$ws.Range("J34").Value2 = 11705.542
$ws.Range("K34").Value2 = 2922.0488
$ws.Range("L34").Value2 = 11705.542
$ws.Range("M34").Value2 = -2720.0488
$ws.Range("N34").Value2 = -12109.542

$ws.Range("H36").Value2 = 49963.332
$ws.Range("I36").Value2 = 0
$ws.Range("J36").Value2 = 49963.332
$ws.Range("K36").Value2 = 0
$ws.Range("L36").Value2 = 49963.332
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value2 = -50739.332

$ws.Range("H39").Value2 = 15000
$ws.Range("I39").Value2 = 0
$ws.Range("J39").Value2 = 15000
$ws.Range("K39").Value2 = 0
$ws.Range("L39").Value2 = 15000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value2 = -15782

$ws.Range("H40").Value2 = 49963.332
$ws.Range("I40").Value2 = 0
$ws.Range("J40").Value2 = 49963.332
$ws.Range("K40").Value2 = 0
$ws.Range("L40").Value2 = 49963.332
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value2 = -50283.332

$ws.Range("H49").Value2 = 15000
$ws.Range("I49").Value2 = 0
$ws.Range("J49").Value2 = 15000
$ws.Range("K49").Value2 = 0
$ws.Range("L49").Value2 = 15000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value2 = -15364

$ws.Range("H134").Value2 = 6834.5
$ws.Range("I134").Value2 = 3104.9333
$ws.Range("J134").Value2 = 9266.825999999999
$ws.Range("K134").Value2 = 9314.7999
$ws.Range("L134").Value2 = 27800.478
$ws.Range("M134").Value2 = -6779.7999
$ws.Range("N134").Value2 = -32870.478

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 222.44444
$ws.Range("I23").Value2 = 153.72728
$ws.Range("K23").Value2 = 461.18184
$ws.Range("M23").Value2 = -226.18184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value2 = 69999.75
$ws.Range("I52").Value2 = 0
$ws.Range("K52").Value2 = 0
$ws.Range("M52").ClearContents()

$ws.Range("H132").Value2 = 8996.8125
$ws.Range("I132").Value2 = 4359.1763
$ws.Range("J132").Value2 = 14252.8
$ws.Range("K132").Value2 = 13077.5289
$ws.Range("L132").Value2 = 42758.39999999999
$ws.Range("M132").Value2 = -10547.5289
$ws.Range("N132").Value2 = -47818.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1861.1052
$ws.Range("J22").Value2 = 4018.5
$ws.Range("L22").Value2 = 4018.5
$ws.Range("N22").Value2 = -4608.5

$ws.Range("H27").Value2 = 1861.1052
$ws.Range("J27").Value2 = 4018.5
$ws.Range("L27").Value2 = 4018.5
$ws.Range("N27").Value2 = -4232.5

$ws.Range("H46").Value2 = 5054440
$ws.Range("J46").Value2 = 5559737
$ws.Range("L46").Value2 = 5559737
$ws.Range("N46").Value2 = -5560113

$ws.Range("H68").Value2 = 5441.4165
$ws.Range("I68").Value2 = 4266.3335
$ws.Range("J68").Value2 = 5833.1113
$ws.Range("K68").Value2 = 4266.3335
$ws.Range("L68").Value2 = 5833.1113
$ws.Range("M68").Value2 = -3517.3335
$ws.Range("N68").Value2 = -7331.1113

$ws.Range("H71").Value2 = 5441.4165
$ws.Range("I71").Value2 = 4266.3335
$ws.Range("J71").Value2 = 5833.1113
$ws.Range("K71").Value2 = 21331.6675
$ws.Range("L71").Value2 = 29165.5565
$ws.Range("M71").Value2 = -17587.6675
$ws.Range("N71").Value2 = -36653.5565

$ws.Range("H82").Value2 = 3574.75
$ws.Range("J82").Value2 = 3340.3635
$ws.Range("L82").Value2 = 3340.3635
$ws.Range("N82").Value2 = -4062.3635

$ws.Range("H85").Value2 = 3574.75
$ws.Range("J85").Value2 = 3340.3635
$ws.Range("L85").Value2 = 3340.3635
$ws.Range("N85").Value2 = -5836.363499999999

$ws.Range("H122").Value2 = 5330.385
$ws.Range("I122").Value2 = 4254.4443
$ws.Range("K122").Value2 = 12763.3329
$ws.Range("M122").Value2 = -10313.3329

$ws.Range("H132").Value2 = 7580948.5
$ws.Range("I132").Value2 = 14288083
$ws.Range("J132").Value2 = 8377.192999999999
$ws.Range("K132").Value2 = 42864249
$ws.Range("L132").Value2 = 25131.579
$ws.Range("M132").Value2 = -42861719
$ws.Range("N132").Value2 = -30191.579

$ws.Range("H136").Value2 = 10195.8
$ws.Range("I136").Value2 = 1942.1428
$ws.Range("K136").Value2 = 5826.428400000001
$ws.Range("M136").Value2 = -3276.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value2 = 22517.5
$ws.Range("J54").Value2 = 25000
$ws.Range("L54").Value2 = 25000
$ws.Range("N54").Value2 = -26040

$ws.Range("H62").Value2 = 6146.5
$ws.Range("I62").Value2 = 6146.5
$ws.Range("K62").Value2 = 6146.5
$ws.Range("M62").Value2 = -5522.5

$ws.Range("H65").Value2 = 6146.5
$ws.Range("I65").Value2 = 6146.5
$ws.Range("K65").Value2 = 30732.5
$ws.Range("M65").Value2 = -27612.5

$ws.Range("H132").Value2 = 22746952
$ws.Range("I132").Value2 = 35724000
$ws.Range("K132").Value2 = 107172000
$ws.Range("M132").Value2 = -107169470
